$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at 297, pushing existing rows 297.. downward.
$ws.Rows("297:298").Insert()

# New row 297
$ws.Cells.Item(297, 1).Value = 3
$ws.Cells.Item(297, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(297, 3).Value = "Coquimbo"
$ws.Cells.Item(297, 4).Value = 44627
$ws.Cells.Item(297, 5).Value = 5
$ws.Cells.Item(297, 6).Value = 100112003
$ws.Cells.Item(297, 7).Value = "Ajo"
$ws.Cells.Item(297, 8).Value = "Chino"
$ws.Cells.Item(297, 9).Value = "1a (cosecha)"
$ws.Cells.Item(297, 10).Value = 69
$ws.Cells.Item(297, 11).Value = 16000
$ws.Cells.Item(297, 12).Value = 16500
$ws.Cells.Item(297, 13).Value = 16254
$ws.Cells.Item(297, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(297, 15).Value = "China"
$ws.Cells.Item(297, 16).Value = 1625
$ws.Cells.Item(297, 17).Value = 10
$ws.Cells.Item(297, 18).Value = "Hortaliza"

# New row 298
$ws.Cells.Item(298, 1).Value = 3
$ws.Cells.Item(298, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(298, 3).Value = "Coquimbo"
$ws.Cells.Item(298, 4).Value = 44627
$ws.Cells.Item(298, 5).Value = 5
$ws.Cells.Item(298, 6).Value = 100112003
$ws.Cells.Item(298, 7).Value = "Ajo"
$ws.Cells.Item(298, 8).Value = "Chino"
$ws.Cells.Item(298, 9).Value = "2a (cosecha)"
$ws.Cells.Item(298, 10).Value = 30
$ws.Cells.Item(298, 11).Value = 15000
$ws.Cells.Item(298, 12).Value = 15000
$ws.Cells.Item(298, 13).Value = 15000
$ws.Cells.Item(298, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(298, 15).Value = "China"
$ws.Cells.Item(298, 16).Value = 1622
$ws.Cells.Item(298, 17).Value = 10
$ws.Cells.Item(298, 18).Value = "Hortaliza"
